$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- BTec logo (header, "first page" header) ---
# Currently displayed/cNvPr name is "image2.jpg" -> should become "image1.jpg"
$hdrFirst = $sec.Headers.Item(2)
$hdrShape = $hdrFirst.Range.InlineShapes.Item(1)
$hdrAsShape = $hdrShape.ConvertToShape()
$hdrAsShape.Name = "image1.jpg"
$hdrAsShape.ConvertToInlineShape() | Out-Null

# --- Pearson logo (default/primary footer) ---
# Currently displayed/cNvPr name is "image1.png" -> should become "image2.png"
$ftrDefault = $sec.Footers.Item(1)
$ftrDefaultShape = $ftrDefault.Range.InlineShapes.Item(1)
$ftrDefaultAsShape = $ftrDefaultShape.ConvertToShape()
$ftrDefaultAsShape.Name = "image2.png"
$ftrDefaultAsShape.ConvertToInlineShape() | Out-Null

# --- Pearson logo (first-page footer) ---
# Currently displayed/cNvPr name is "image1.png" -> should become "image2.png"
$ftrFirst = $sec.Footers.Item(2)
$ftrFirstShape = $ftrFirst.Range.InlineShapes.Item(1)
$ftrFirstAsShape = $ftrFirstShape.ConvertToShape()
$ftrFirstAsShape.Name = "image2.png"
$ftrFirstAsShape.ConvertToInlineShape() | Out-Null

Write-Output "Renamed BTec logo and both Pearson logo inline pictures."
